$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16128.74827949199
$ws.Range("C2").Value = 3145.064652637094
$ws.Range("D2").Value = 9689.736623685647
